$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with updated crypto market data.
# Price values that look like plain numbers (e.g. "59.33") are written with a leading
# apostrophe (quote-prefix) so Excel stores them as text -- matching the original
# workbook, where every Price/Volume cell is text -- instead of coercing them to Double
# (which would also lose precision, e.g. 59.33 -> 59.329999999999998).

$ws.Range("D2").Value = '37.167.37'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").Value = '2.077.08'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "`'253.47"
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("E6").Value = '  +1.93%  '
$ws.Range("D7").Value = "`'59.33"
$ws.Range("E7").Value = '  +9.13%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").Value = "`'0.393"
$ws.Range("E9").Value = '  +4.79%  '
$ws.Range("D10").Value = "`'61.54"
$ws.Range("E10").Value = '  -0.42%  '
$ws.Range("D11").Value = "`'0.0803"
$ws.Range("E11").Value = '  +8.03%  '
$ws.Range("E12").Value = '  +2.55%  '
$ws.Range("D13").Value = "`'16.35"
$ws.Range("E13").Value = '  +7.46%  '
$ws.Range("D14").Value = '2.379.41'
$ws.Range("E14").Value = '  -0.89%  '
$ws.Range("D15").Value = "`'0.818"
$ws.Range("E15").Value = '  -2.14%  '
$ws.Range("D16").Value = "`'5.58"
$ws.Range("E16").Value = '  +7.66%  '
$ws.Range("D17").Value = '2.075.98'
$ws.Range("E17").Value = '  -0.99%  '
$ws.Range("D18").Value = '37.274.40'
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").Value = "`'15.96"
$ws.Range("E19").Value = '  +8.82%  '
$ws.Range("D20").Value = "`'74.93"
$ws.Range("E20").Value = '  +2.88%  '
$ws.Range("D21").Value = '0.0₃0933'
$ws.Range("E21").Value = '  +10.25%  '
$ws.Range("D22").Value = "`'5.49"
$ws.Range("E22").Value = '  +5.37%  '
$ws.Range("D23").Value = "`'239.44"
$ws.Range("E23").Value = '  -0.73%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").Value = "`'2.42"
$ws.Range("E25").Value = '  -2.31%  '
$ws.Range("D26").Value = "`'2.29"
$ws.Range("E26").Value = '  +14.49%  '
$ws.Range("D27").Value = "`'170.41"
$ws.Range("E27").Value = '  -1.10%  '
$ws.Range("D28").Value = "`'9.39"
$ws.Range("E28").Value = '  +1.41%  '
$ws.Range("D29").Value = "`'20.44"
$ws.Range("E29").Value = '  -1.16%  '
$ws.Range("E30").Value = '  +3.10%  '
$ws.Range("E31").Value = '  +6.86%  '
$ws.Range("E32").Value = '  +6.07%  '
$ws.Range("D33").Value = "`'0.0638"
$ws.Range("E33").Value = '  +3.46%  '
$ws.Range("E34").Value = '  +9.06%  '
$ws.Range("D35").Value = "`'0.0912"
$ws.Range("E35").Value = '  +0.22%  '
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("E37").Value = '  +2.29%  '
$ws.Range("D38").Value = "`'0.118"
$ws.Range("E38").Value = '  +26.93%  '
$ws.Range("E39").Value = '  -3.84%  '
$ws.Range("E40").Value = '  +2.49%  '
$ws.Range("D41").Value = "`'0.0228"
$ws.Range("E41").Value = '  +0.56%  '
$ws.Range("E42").Value = '  -2.57%  '
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").Value = "`'99.32"
$ws.Range("E44").Value = '  +0.49%  '
$ws.Range("D45").Value = "`'4.36"
$ws.Range("E45").Value = '  +2.90%  '
$ws.Range("D46").Value = "`'2.85"
$ws.Range("E46").Value = '  +1.50%  '
$ws.Range("D47").Value = "`'4.57"
$ws.Range("E47").Value = '  +13.34%  '
$ws.Range("E48").Value = '  +7.91%  '
$ws.Range("D49").Value = '1.309.97'
$ws.Range("E49").Value = '  -0.98%  '
$ws.Range("E50").Value = '  -0.11%  '
$ws.Range("D51").Value = "`'6.95"
$ws.Range("E51").Value = '  -0.73%  '
